# Auto update Excel log
# Appends new sensor-log rows (the "2026-02-01 19:55/19:56" batch of
# readings) to the PIR, Humidity, and Temperature sheets, extending each
# sheet's used range accordingly.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($Sheet, $Rows)

    foreach ($row in $Rows) {
        $r = $row[0]

        # Column A ("Date") looks like a date string (e.g. "2026-02-01").
        # Force text formatting first so Excel stores the literal string
        # instead of silently converting it to a date serial number.
        $aCell = $Sheet.Cells.Item($r, 1)
        $aCell.NumberFormat = "@"
        $aCell.Value = $row[1]

        $Sheet.Cells.Item($r, 2).Value = $row[2]
        $Sheet.Cells.Item($r, 3).Value = $row[3]
        $Sheet.Cells.Item($r, 4).Value = $row[4]

        # Column E ("Value") can look like a percentage (e.g. "78.4%"),
        # which Excel would otherwise auto-convert to a numeric/percent
        # cell. Force text formatting there too so the literal string is
        # preserved (a no-op for non-numeric values like "No Motion").
        $eCell = $Sheet.Cells.Item($r, 5)
        $eCell.NumberFormat = "@"
        $eCell.Value = $row[5]

        $Sheet.Cells.Item($r, 6).Value = $row[6]
    }
}

# ---------------------------------------------------------------------------
# PIR sheet: append rows 43-55 (A1:F42 -> A1:F55)
# ---------------------------------------------------------------------------
$pirRows = @(
    @(43, "2026-02-01", "19:55:58", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(44, "2026-02-01", "19:56:01", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(45, "2026-02-01", "19:56:06", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(46, "2026-02-01", "19:56:11", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(47, "2026-02-01", "19:56:16", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(48, "2026-02-01", "19:56:21", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(49, "2026-02-01", "19:56:27", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(50, "2026-02-01", "19:56:32", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(51, "2026-02-01", "19:56:36", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(52, "2026-02-01", "19:56:42", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(53, "2026-02-01", "19:56:47", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(54, "2026-02-01", "19:56:52", "19:00", "Bathroom", "No Motion", "Inactive"),
    @(55, "2026-02-01", "19:56:57", "19:00", "Bathroom", "No Motion", "Inactive")
)
$wsPIR = $wb.Worksheets.Item("PIR")
Add-LogRows $wsPIR $pirRows

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 35-44 (A1:F34 -> A1:F44)
# ---------------------------------------------------------------------------
$humidityRows = @(
    @(35, "2026-02-01", "19:55:58", "19:00", "Bathroom", "78.4%", "Active"),
    @(36, "2026-02-01", "19:56:01", "19:00", "Bathroom", "77.8%", "Active"),
    @(37, "2026-02-01", "19:56:11", "19:00", "Bathroom", "78.0%", "Active"),
    @(38, "2026-02-01", "19:56:21", "19:00", "Bathroom", "77.9%", "Active"),
    @(39, "2026-02-01", "19:56:26", "19:00", "Bathroom", "79.0%", "Active"),
    @(40, "2026-02-01", "19:56:31", "19:00", "Bathroom", "77.6%", "Active"),
    @(41, "2026-02-01", "19:56:36", "19:00", "Bathroom", "78.9%", "Active"),
    @(42, "2026-02-01", "19:56:41", "19:00", "Bathroom", "78.5%", "Active"),
    @(43, "2026-02-01", "19:56:46", "19:00", "Bathroom", "79.3%", "Active"),
    @(44, "2026-02-01", "19:56:56", "19:00", "Bathroom", "78.7%", "Active")
)
$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRows $wsHumidity $humidityRows

# ---------------------------------------------------------------------------
# Temperature sheet: append rows 35-44 (A1:F34 -> A1:F44)
# ---------------------------------------------------------------------------
$temperatureRows = @(
    @(35, "2026-02-01", "19:55:58", "19:00", "Bathroom", "25.2C", "Active"),
    @(36, "2026-02-01", "19:56:01", "19:00", "Bathroom", "25.2C", "Active"),
    @(37, "2026-02-01", "19:56:11", "19:00", "Bathroom", "25.2C", "Active"),
    @(38, "2026-02-01", "19:56:21", "19:00", "Bathroom", "25.2C", "Active"),
    @(39, "2026-02-01", "19:56:26", "19:00", "Bathroom", "25.2C", "Active"),
    @(40, "2026-02-01", "19:56:31", "19:00", "Bathroom", "25.2C", "Active"),
    @(41, "2026-02-01", "19:56:36", "19:00", "Bathroom", "25.2C", "Active"),
    @(42, "2026-02-01", "19:56:42", "19:00", "Bathroom", "25.2C", "Active"),
    @(43, "2026-02-01", "19:56:46", "19:00", "Bathroom", "25.2C", "Active"),
    @(44, "2026-02-01", "19:56:56", "19:00", "Bathroom", "25.2C", "Active")
)
$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRows $wsTemperature $temperatureRows
